$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row with no label (A column empty) holding "519033 - Carlos Yujiro Shigue"
# under "Docentes responsaveis:" (row 13) is removed; everything below shifts up.
$ws.Rows.Item(13).Delete()

# After the shift, update the contents of a few description cells.
# Objetivos: (row 10) now holds the professor identification string.
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

# Programa resumido: (row 13, after shift) now holds "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Programa: (row 15, after shift) now holds "01/01/2012". Assigning that text
# directly would be auto-parsed as a date by Excel, so copy it from a cell
# that already holds it as plain text to keep it a shared string.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

# Metodo: (row 18, after shift) now holds the professor identification string.
$ws.Range("B18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C18").Value = "519033 - Carlos Yujiro Shigue"
